$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their original "text" storage type
# (values like "1.004" or "0.000008190" must remain literal text, not be
# reinterpreted as numbers, which would lose formatting / precision).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.550.62'
$ws.Range("E2").Value = '  -7.53%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.685.25'
$ws.Range("E3").Value = '  -6.64%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.15%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.01'
$ws.Range("E5").Value = '  -6.36%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.4995'
$ws.Range("E6").Value = '  -15.93%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.004'
$ws.Range("E7").Value = '  +0.03%  '

# Row 8
$ws.Range("E8").Value = '  -6.38%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '21.93'
$ws.Range("E9").Value = '  -6.02%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06196'
$ws.Range("E10").Value = '  -9.45%  '

# Row 11
$ws.Range("E11").Value = '  -3.35%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.676.37'
$ws.Range("E12").Value = '  -6.80%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.444'
$ws.Range("E13").Value = '  -7.67%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5757'
$ws.Range("E14").Value = '  -7.70%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.914.33'
$ws.Range("E15").Value = '  -6.64%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008190'
$ws.Range("E16").Value = '  -12.33%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.66'
$ws.Range("E17").Value = '  -14.67%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.581.17'
$ws.Range("E18").Value = '  -7.26%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.000'
$ws.Range("E19").Value = '  -9.04%  '

# Row 20
$ws.Range("E20").Value = '  +0.09%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.77'
$ws.Range("E21").Value = '  -6.16%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '185.06'
$ws.Range("E22").Value = '  -12.18%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.195'
$ws.Range("E23").Value = '  -9.84%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.005'
$ws.Range("E24").Value = '  +0.10%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.54'
$ws.Range("E25").Value = '  -6.31%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.488'
$ws.Range("E26").Value = '  -4.98%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1133'
$ws.Range("E27").Value = '  -11.16%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.46'
$ws.Range("E28").Value = '  -5.86%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.299'
$ws.Range("E29").Value = '  -8.85%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05699'
$ws.Range("E30").Value = '  -8.19%  '

# Row 31
$ws.Range("E31").Value = '  -7.32%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.478'
$ws.Range("E32").Value = '  -8.15%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.477'
$ws.Range("E33").Value = '  -7.50%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.637'
$ws.Range("E34").Value = '  -4.97%  '

# Row 35
$ws.Range("E35").Value = '  -5.56%  '

# Row 36
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.369'
$ws.Range("E36").Value = '  -5.00%  '

# Row 37
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5919'
$ws.Range("E37").Value = '  -7.92%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.635'
$ws.Range("E38").Value = '  -2.95%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01589'
$ws.Range("E39").Value = '  -7.36%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.069.29'
$ws.Range("E40").Value = '  -5.86%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.888'
$ws.Range("E41").Value = '  -9.32%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8546'
$ws.Range("E42").Value = '  -2.88%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  -0.63%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '98.15'
$ws.Range("E44").Value = '  -2.47%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.839.85'
$ws.Range("E45").Value = '  -6.28%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.34'
$ws.Range("E46").Value = '  -6.99%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000105'
$ws.Range("E47").Value = '  -6.37%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.005'
$ws.Range("E48").Value = '  -0.07%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.007'
$ws.Range("E49").Value = '  -4.14%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4312'
$ws.Range("E50").Value = '  -3.83%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05196'
$ws.Range("E51").Value = '  -4.93%  '
